$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(1.02, 1.078283938229475, 1.074571689955286, 1.080817204628942, 1.08924329764285, 1.063968640316946, 1.083174181525269, 1.077260549269478, 1.083489648568651, 1.091893916778798, 1.08471241310485)
    3 = @(1.02, 1.079645288707404, 1.075445756466667, 1.081981253421026, 1.090358272207764, 1.064412962390437, 1.084194143277961, 1.077951626299769, 1.084471187776625, 1.092828080983467, 1.085733823320222)
    4 = @(1.02, 1.080525682286966, 1.076010897528084, 1.082734143397934, 1.091079497774404, 1.064698939183223, 1.084853115502467, 1.078397718892392, 1.085105396899392, 1.093431713187621, 1.086393731360872)
    5 = @(1.02, 1.08089568566833, 1.076248378854061, 1.083050582162135, 1.091382645906755, 1.064818798737808, 1.085129907685734, 1.0785849986946, 1.0853718020588, 1.0936852820943, 1.086670916620788)
    6 = @(1.02, 1.080957804276632, 1.076288246944889, 1.083103709219707, 1.091433542632605, 1.064838902305281, 1.085176368321043, 1.078616428729209, 1.085416520030178, 1.093727845864064, 1.086717443235528)
    7 = @(1.02, 1.080530626734124, 1.076014071172525, 1.082738371964786, 1.091083548668873, 1.064700542185731, 1.084856814952256, 1.078400222344754, 1.085108957465727, 1.093435102163925, 1.086397436064304)
    8 = @(1.02, 1.078744116824203, 1.074867176195765, 1.081210669493478, 1.089620158767348, 1.064119118313351, 1.083519093174008, 1.077494326390703, 1.083821554500429, 1.092209795765199, 1.085057814567692)
    9 = @(1.02, 1.07559213835523, 1.072842798356619, 1.078516055709397, 1.087039579805649, 1.063082816265241, 1.081154006062022, 1.075889695506703, 1.081545905270392, 1.090044181075678, 1.082689368760728)
    10 = @(1.02, 1.073487961692801, 1.071490860639548, 1.076717752956773, 1.085317803002058, 1.062383976172358, 1.079571853546805, 1.074814269410431, 1.080023913091883, 1.088595979238511, 1.081104969407465)
    11 = @(1.02, 1.072576102132534, 1.070904883444002, 1.07593858682275, 1.084571898578918, 1.062079464706291, 1.078885442737491, 1.074347235598203, 1.079363683184279, 1.08796780917987, 1.080417583816069)
    12 = @(1.02, 1.072237281741689, 1.070687136987171, 1.075649093479795, 1.08429477946964, 1.061966067397612, 1.078630276622055, 1.074173551480362, 1.07911826185083, 1.087734312966962, 1.08016205533548)
    13 = @(1.02, 1.072309965147198, 1.070733848365991, 1.075711194285467, 1.084354225076332, 1.06199040456863, 1.078685019835581, 1.074210816703917, 1.079170913863756, 1.08778440624021, 1.080216876290648)
    14 = @(1.02, 1.07254809750783, 1.070886886266214, 1.075914658786753, 1.084548993003227, 1.062070097136359, 1.078864354776787, 1.074332883045901, 1.079343400319751, 1.087948511707146, 1.080396465908042)
    15 = @(1.02, 1.072694803324366, 1.07098116620645, 1.076040009778658, 1.084668988328167, 1.062119160128189, 1.078974822018973, 1.074408064676291, 1.079449650648338, 1.088049600452557, 1.080507090026382)
    16 = @(1.02, 1.073548462902199, 1.071529737691953, 1.076769453097148, 1.085367298358654, 1.062404145256061, 1.079617380128374, 1.074845235968058, 1.080067704923586, 1.088637645677755, 1.081150560642)
    17 = @(1.02, 1.074083740538007, 1.071873686018496, 1.077226880640316, 1.085805230759645, 1.062582396811504, 1.080020082099336, 1.075119094830085, 1.080455071617332, 1.089006217735434, 1.081553834495953)
    18 = @(1.02, 1.074395888098356, 1.072074249466588, 1.077493643383414, 1.086060634388256, 1.062686183821598, 1.08025484338623, 1.075278700260843, 1.080680900720922, 1.089221094651951, 1.081788929170801)
    19 = @(1.02, 1.074502310468443, 1.072142626998265, 1.077584594698158, 1.086147714570923, 1.062721541323503, 1.080334869304679, 1.075333099277517, 1.080757883126864, 1.089294344397555, 1.081869068735234)
    20 = @(1.02, 1.074026317672894, 1.071836789399997, 1.077177807877642, 1.085758248388826, 1.062563291160075, 1.079976889247376, 1.075089726011734, 1.080413522787182, 1.088966684312404, 1.081510580305188)
    21 = @(1.02, 1.072477976676994, 1.07084182287161, 1.075854745692923, 1.084491640295515, 1.062046637627996, 1.078811550722809, 1.074296943276405, 1.079292612418358, 1.087900191334744, 1.080343586866251)
    22 = @(1.02, 1.071503804831599, 1.070215735557026, 1.07502244083471, 1.083694942586334, 1.061720128974639, 1.078077682584572, 1.073797291091529, 1.078586794076239, 1.087228683992642, 1.07960867655108)
    23 = @(1.02, 1.072020296586988, 1.070547685328186, 1.07546370434889, 1.084117319271977, 1.061893375962433, 1.078466832217219, 1.074062280266478, 1.078961062868363, 1.087584754575978, 1.079998378820842)
    24 = @(1.02, 1.07405226481257, 1.071853461566149, 1.077199981891973, 1.085779477807805, 1.062571924750398, 1.079996406633247, 1.075102996923861, 1.080432297271386, 1.088984548088769, 1.081530125407986)
    25 = @(1.02, 1.076407487698094, 1.073366558401841, 1.079213001357134, 1.087706958071534, 1.063352125459274, 1.0817663834568, 1.076305525369409, 1.082135067808673, 1.090604822181247, 1.083302615801647)
}

$cols = @("B","C","D","E","F","I","J","K","L","M","N")

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = "$($cols[$i])$row"
        $ws.Range($addr).Value = $vals[$i]
    }
}
